$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "[-, -, 'ELM-2NA-Eletrohidráulica', -]"
$ws.Range("D18").Value = "[-, -, -, 'ELM-2NA-Eletrohidráulica']"
$ws.Range("F18").Value = "-"

$ws.Range("B19").Value = "[-, -, 'ELM-2NA-Eletrohidráulica', -]"
$ws.Range("D19").Value = "[-, -, -, 'ELM-2NA-Eletropneumática']"
$ws.Range("F19").Value = "-"

$ws.Range("B20").Value = "[-, 'ELM-2NA-Eletropneumática', -, -]"
$ws.Range("D20").Value = "[-, -, -, 'ELM-2NA-Eletropneumática']"
$ws.Range("F20").Value = "-"

$ws.Range("B21").Value = "['ELM-2NA-Eletrohidráulica', -, -, -]"
$ws.Range("D21").Value = "[-, -, -, 'ELM-2NA-Eletropneumática']"
$ws.Range("F21").Value = "-"
